{"js": "// Replace the three-digit-by-one-digit multiplication equations in the\n// document's table cells with their updated answers. Each original\n// equation string is unique within the document, so a direct\n// search-and-replace mapping is safe and unambiguous.\nconst replacements = [\n  [\"633\u00d74=2532\", \"125\u00d75=625\"],\n  [\"476\u00d75=2380\", \"760\u00d78=6080\"],\n  [\"640\u00d78=5120\", \"788\u00d78=6304\"],\n  [\"887\u00d77=6209\", \"888\u00d77=6216\"],\n  [\"266\u00d79=2394\", \"102\u00d77=714\"],\n  [\"846\u00d77=5922\", \"874\u00d76=5244\"],\n  [\"305\u00d79=2745\", \"887\u00d73=2661\"],\n  [\"781\u00d73=2343\", \"770\u00d72=1540\"],\n  [\"666\u00d76=3996\", \"925\u00d76=5550\"],\n  [\"490\u00d79=4410\", \"300\u00d77=2100\"],\n  [\"880\u00d73=2640\", \"759\u00d74=3036\"],\n  [\"552\u00d77=3864\", \"718\u00d75=3590\"],\n  [\"193\u00d72=386\", \"467\u00d72=934\"],\n  [\"499\u00d73=1497\", \"711\u00d72=1422\"],\n  [\"368\u00d73=1104\", \"738\u00d79=6642\"],\n  [\"346\u00d76=2076\", \"127\u00d79=1143\"],\n  [\"139\u00d74=556\", \"317\u00d77=2219\"],\n  [\"233\u00d77=1631\", \"779\u00d73=2337\"],\n  [\"690\u00d78=5520\", \"265\u00d79=2385\"],\n  [\"816\u00d77=5712\", \"333\u00d75=1665\"],\n  [\"788\u00d75=3940\", \"985\u00d77=6895\"],\n  [\"972\u00d77=6804\", \"359\u00d75=1795\"],\n  [\"321\u00d76=1926\", \"714\u00d76=4284\"],\n  [\"364\u00d78=2912\", \"318\u00d78=2544\"],\n  [\"384\u00d73=1152\", \"404\u00d77=2828\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication equations in the\n# document's table cells with their updated answers. Each original\n# equation string is unique within the document, so a direct\n# search-and-replace mapping is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"633\u00d74=2532\", \"125\u00d75=625\"),\n    @(\"476\u00d75=2380\", \"760\u00d78=6080\"),\n    @(\"640\u00d78=5120\", \"788\u00d78=6304\"),\n    @(\"887\u00d77=6209\", \"888\u00d77=6216\"),\n    @(\"266\u00d79=2394\", \"102\u00d77=714\"),\n    @(\"846\u00d77=5922\", \"874\u00d76=5244\"),\n    @(\"305\u00d79=2745\", \"887\u00d73=2661\"),\n    @(\"781\u00d73=2343\", \"770\u00d72=1540\"),\n    @(\"666\u00d76=3996\", \"925\u00d76=5550\"),\n    @(\"490\u00d79=4410\", \"300\u00d77=2100\"),\n    @(\"880\u00d73=2640\", \"759\u00d74=3036\"),\n    @(\"552\u00d77=3864\", \"718\u00d75=3590\"),\n    @(\"193\u00d72=386\",  \"467\u00d72=934\"),\n    @(\"499\u00d73=1497\", \"711\u00d72=1422\"),\n    @(\"368\u00d73=1104\", \"738\u00d79=6642\"),\n    @(\"346\u00d76=2076\", \"127\u00d79=1143\"),\n    @(\"139\u00d74=556\",  \"317\u00d77=2219\"),\n    @(\"233\u00d77=1631\", \"779\u00d73=2337\"),\n    @(\"690\u00d78=5520\", \"265\u00d79=2385\"),\n    @(\"816\u00d77=5712\", \"333\u00d75=1665\"),\n    @(\"788\u00d75=3940\", \"985\u00d77=6895\"),\n    @(\"972\u00d77=6804\", \"359\u00d75=1795\"),\n    @(\"321\u00d76=1926\", \"714\u00d76=4284\"),\n    @(\"364\u00d78=2912\", \"318\u00d78=2544\"),\n    @(\"384\u00d73=1152\", \"404\u00d77=2828\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
